$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 277.0476
$ws.Range("I53").Value = 75.36364
$ws.Range("J53").Value = 498.9
$ws.Range("K53").Value = 75.36364
$ws.Range("L53").Value = 498.9
$ws.Range("M53").Value = 561.63636
$ws.Range("N53").Value = -1772.9
$ws.Range("H64").Value = 4322.6313
$ws.Range("J64").Value = 4471.5386
$ws.Range("L64").Value = 4471.5386
$ws.Range("N64").Value = -4967.5386
$ws.Range("H67").Value = 4322.6313
$ws.Range("J67").Value = 4471.5386
$ws.Range("L67").Value = 4471.5386
$ws.Range("N67").Value = -6187.5386
$ws.Range("H74").Value = 5618.421
$ws.Range("I74").Value = 14232.667
$ws.Range("J74").Value = 4003.25
$ws.Range("K74").Value = 14232.667
$ws.Range("L74").Value = 4003.25
$ws.Range("M74").Value = -13296.667
$ws.Range("N74").Value = -5875.25
$ws.Range("H77").Value = 5618.421
$ws.Range("I77").Value = 14232.667
$ws.Range("J77").Value = 4003.25
$ws.Range("K77").Value = 71163.33499999999
$ws.Range("L77").Value = 20016.25
$ws.Range("M77").Value = -66483.33499999999
$ws.Range("N77").Value = -29376.25
$ws.Range("H129").Value = 938.3929000000001
$ws.Range("I129").Value = 776.5
$ws.Range("J129").Value = 1154.25
$ws.Range("K129").Value = 2329.5
$ws.Range("L129").Value = 3462.75
$ws.Range("M129").Value = 2670.5
$ws.Range("N129").Value = -13462.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 683.2778
$ws.Range("I45").Value = 542.7857
$ws.Range("J45").Value = 1175
$ws.Range("K45").Value = 542.7857
$ws.Range("L45").Value = 1175
$ws.Range("M45").Value = -165.7857
$ws.Range("N45").Value = -1929
$ws.Range("H61").Value = 35788100
$ws.Range("I61").Value = 45501890
$ws.Range("J61").Value = 170883.33
$ws.Range("K61").Value = 45501890
$ws.Range("L61").Value = 170883.33
$ws.Range("M61").Value = -45501678
$ws.Range("N61").Value = -171307.33
$ws.Range("H74").Value = 8599917
$ws.Range("I74").Value = 11145532
$ws.Range("J74").Value = 114533
$ws.Range("K74").Value = 11145532
$ws.Range("L74").Value = 114533
$ws.Range("M74").Value = -11144658
$ws.Range("N74").Value = -116281
$ws.Range("H77").Value = 8599917
$ws.Range("I77").Value = 11145532
$ws.Range("J77").Value = 114533
$ws.Range("K77").Value = 55727660
$ws.Range("L77").Value = 572665
$ws.Range("M77").Value = -55723292
$ws.Range("N77").Value = -581401
$ws.Range("H136").Value = 35788100
$ws.Range("I136").Value = 45501890
$ws.Range("J136").Value = 170883.33
$ws.Range("K136").Value = 136505670
$ws.Range("L136").Value = 512649.99
$ws.Range("M136").Value = -136503120
$ws.Range("N136").Value = -517749.99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2845.4783
$ws.Range("I134").Value = 2804.1052
$ws.Range("K134").Value = 8412.3156
$ws.Range("M134").Value = -5877.3156
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4341.0835
$ws.Range("I62").Value = 2978.6
$ws.Range("J62").Value = 5314.2856
$ws.Range("K62").Value = 2978.6
$ws.Range("L62").Value = 5314.2856
$ws.Range("M62").Value = -2354.6
$ws.Range("N62").Value = -6562.2856
$ws.Range("H65").Value = 4341.0835
$ws.Range("I65").Value = 2978.6
$ws.Range("J65").Value = 5314.2856
$ws.Range("K65").Value = 14893
$ws.Range("L65").Value = 26571.428
$ws.Range("M65").Value = -11773
$ws.Range("N65").Value = -32811.428
$ws.Range("H105").Value = 836.2727
$ws.Range("I105").Value = 836.2727
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 836.2727
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 910.7273
$ws.Range("N105").ClearContents()
$ws.Range("H134").Value = 17093.338
$ws.Range("I134").Value = 1895.8163
$ws.Range("J134").Value = 50942.363
$ws.Range("K134").Value = 5687.448899999999
$ws.Range("L134").Value = 152827.089
$ws.Range("M134").Value = -3152.448899999999
$ws.Range("N134").Value = -157897.089
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 681.2162
$ws.Range("I68").Value = 691.86664
$ws.Range("J68").Value = 635.5714
$ws.Range("K68").Value = 2075.59992
$ws.Range("L68").Value = 1906.7142
$ws.Range("M68").Value = -1264.59992
$ws.Range("N68").Value = -3528.7142
$ws.Range("H71").Value = 681.2162
$ws.Range("I71").Value = 691.86664
$ws.Range("J71").Value = 635.5714
$ws.Range("K71").Value = 6226.79976
$ws.Range("L71").Value = 5720.1426
$ws.Range("M71").Value = -2170.79976
$ws.Range("N71").Value = -13832.1426
$ws.Range("H122").Value = 1003.6
$ws.Range("I122").Value = 355.33334
$ws.Range("J122").Value = 1368.25
$ws.Range("K122").Value = 3198.00006
$ws.Range("L122").Value = 12314.25
$ws.Range("M122").Value = -748.0000600000003
$ws.Range("N122").Value = -17214.25
$ws.Range("H131").Value = 904.0714
$ws.Range("J131").Value = 1041.5454
$ws.Range("L131").Value = 3124.6362
$ws.Range("N131").Value = -13204.6362
$ws.Range("H132").Value = 2446.7144
$ws.Range("I132").Value = 1269.5385
$ws.Range("J132").Value = 3466.9333
$ws.Range("K132").Value = 11425.8465
$ws.Range("L132").Value = 31202.3997
$ws.Range("M132").Value = -8895.846500000001
$ws.Range("N132").Value = -36262.3997
$ws.Range("H140").Value = 2477.558
$ws.Range("I140").Value = 3112.8096
$ws.Range("J140").Value = 1871.1818
$ws.Range("K140").Value = 9338.4288
$ws.Range("L140").Value = 5613.5454
$ws.Range("M140").Value = -4158.4288
$ws.Range("N140").Value = -15973.5454
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 82290.75999999999
$ws.Range("I132").Value = 51913.8
$ws.Range("J132").Value = 203798.6
$ws.Range("K132").Value = 155741.4
$ws.Range("L132").Value = 611395.8
$ws.Range("M132").Value = -153211.4
$ws.Range("N132").Value = -616455.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1611.375
$ws.Range("I82").Value = 1130.3334
$ws.Range("K82").Value = 1130.3334
$ws.Range("M82").Value = -769.3334
$ws.Range("H85").Value = 1611.375
$ws.Range("I85").Value = 1130.3334
$ws.Range("K85").Value = 1130.3334
$ws.Range("M85").Value = 117.6666
$ws.Range("H93").Value = 1507.0769
$ws.Range("I93").Value = 1424.3334
$ws.Range("K93").Value = 1424.3334
$ws.Range("M93").Value = -176.3334
$ws.Range("H100").Value = 1292.0741
$ws.Range("I100").Value = 1005.0526
$ws.Range("K100").Value = 1005.0526
$ws.Range("M100").Value = -464.0526
$ws.Range("H136").Value = 50295.254
$ws.Range("I136").Value = 31322.111
$ws.Range("J136").Value = 147871.42
$ws.Range("K136").Value = 93966.333
$ws.Range("L136").Value = 443614.26
$ws.Range("M136").Value = -91416.333
$ws.Range("N136").Value = -448714.26
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 52936.41
$ws.Range("I136").Value = 38193.332
$ws.Range("J136").Value = 86108.336
$ws.Range("K136").Value = 114579.996
$ws.Range("L136").Value = 258325.008
$ws.Range("M136").Value = -112029.996
$ws.Range("N136").Value = -263425.008
